$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header + quantity values (A1:A4)
$ws.Range("A1").Value = "Quantity"
$ws.Range("A2").Value = 2
$ws.Range("A3").Value = 3
$ws.Range("A4").Value = 4

# Give A2 a thin black border on the right/top/bottom (not left) plus the
# "unlocked" protection that already exists on the sheet's other style.
$ws.Range("A2").Borders.Item(10).Color = 0   # xlEdgeRight
$ws.Range("A2").Borders.Item(8).Color = 0    # xlEdgeTop
$ws.Range("A2").Borders.Item(9).Color = 0    # xlEdgeBottom

# Propagate the same format to A3:A4 via copy/paste-special so every cell
# lands on the very same style record instead of building its own chain of
# one-off styles.
$ws.Range("A2").Copy()
$ws.Range("A3:A4").PasteSpecial(-4122)  # xlPasteFormats

# Column A is very slightly wider than the sheet default.
$ws.Columns("A").ColumnWidth = 8.64

# Put the live selection on A2, then restore Action1 as the active tab so
# only Global's sheetView selection element actually changes.
$ws.Range("A2").Select()
$wb.Worksheets.Item(2).Select()
